# Sakurazaka_database.xlsx update
# - mmb_info: move "Center" (SP Post) flag from MORITA HIKARU (row16) to YAMAZAKI TEN (row18)
# - song_info: fill in Center for "Time Machine de Yeah!" (row31 -> SUGAI YUUKA) and
#   "Zutto Haru Dattara Naa" (row32 -> TAMURA HONO); add two new songs (rows 55 & 56)
# - restore the selection / active-sheet state to match the saved workbook

$wb = $excel.ActiveWorkbook

$mmb  = $wb.Worksheets.Item("mmb_info")
$song = $wb.Worksheets.Item("song_info")

# --- mmb_info: Center moves from MORITA HIKARU to YAMAZAKI TEN -------------
$mmb.Cells.Item(16, 3).ClearContents()       # C16 (MORITA HIKARU)  was "Center"
$mmb.Cells.Item(18, 3).Value = "Center"      # C18 (YAMAZAKI TEN)   now "Center"

# --- song_info: fill in the two previously-blank ("null") centers ----------
$song.Cells.Item(31, 5).Value = "SUGAI YUUKA"   # Time Machine de Yeah!
$song.Cells.Item(32, 5).Value = "TAMURA HONO"   # Zutto Haru Dattara Naa

# --- song_info: two new rows appended ---------------------------------------
# (cell-write order below matches the original authoring order so new shared
#  strings land at the same indices as the source workbook)
$song.Cells.Item(55, 2).Value = "Kimiga Sayonarawo Ietatte"
$song.Cells.Item(55, 1).Value = "君がサヨナラ言えたって・・・"
$song.Cells.Item(55, 3).Value = "others"
$song.Cells.Item(55, 4).Value = "7th single"
$song.Cells.Item(55, 5).Value = "KOBAYASHI YUI"

$song.Cells.Item(56, 4).Value = "8th single"
$song.Cells.Item(56, 1).Value = "何歳の頃に戻りたいのか？"
$song.Cells.Item(56, 2).Value = "What age do you want to go back to?"
$song.Cells.Item(56, 3).Value = "TITLE"
$song.Cells.Item(56, 5).Value = "YAMAZAKI TEN"

# --- restore view/selection state -------------------------------------------
$song.Activate()
$song.Range("E57").Select()

$mmb.Activate()
$mmb.Range("C13").Select()
